# Database wilayah Indonesia 2018 Semester 1
# Rename the "bps_name" column header (B2) to "name", and leave the
# worksheet scrolled/selected on the header cell B2 (matching the
# saved view state in the target workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B2 currently holds "bps_name" -> rename to "name".
$ws.Range("B2").Value = "name"

# Bring the view back to the top and select the renamed header cell.
$ws.Range("B2").Select()
